$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2035556666666667
$ws.Range("H2").Value = 0.6106670000000001
$ws.Range("I2").Value = 0.006148914270823412
$ws.Range("J2").Value = 0.006148914270823412
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.004739997254000001
$ws.Range("R2").Value = 0.042659975286
$ws.Range("S2").Value = [double]"5.725076051612372e-05"
$ws.Range("T2").Value = [double]"5.725076051612373e-05"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2035556666666667
$ws.Range("H3").Value = 0.6106670000000001
$ws.Range("I3").Value = 0.006148914270823412
$ws.Range("J3").Value = 0.006148914270823412
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.02705648350955555
$ws.Range("R3").Value = 0.243508351586
$ws.Range("S3").Value = 0.0003267943365382412
$ws.Range("T3").Value = 0.0003267943365382413
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2035556666666667
$ws.Range("H4").Value = 0.6106670000000001
$ws.Range("I4").Value = 0.006148914270823412
$ws.Range("J4").Value = 0.006148914270823412
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.477294341716889
$ws.Range("R4").Value = 4.295649075452001
$ws.Range("S4").Value = 0.005764869173769047
$ws.Range("T4").Value = 0.005764869173769047
$ws.Range("I5").Value = 0.735846381812327
$ws.Range("J5").Value = 0.735846381812327
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.5672399509139999
$ws.Range("R5").Value = 5.105159558225999
$ws.Range("S5").Value = 0.006851252615716216
$ws.Range("T5").Value = 0.006851252615716217
$ws.Range("I6").Value = 0.735846381812327
$ws.Range("J6").Value = 0.735846381812327
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.03910778709006508
$ws.Range("T6").Value = 0.03910778709006509
$ws.Range("I7").Value = 0.735846381812327
$ws.Range("J7").Value = 0.735846381812327
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 57.11826493961467
$ws.Range("R7").Value = 514.064384456532
$ws.Range("S7").Value = 0.6898873421065457
$ws.Range("T7").Value = 0.6898873421065457
$ws.Range("G8").Value = 8.541072
$ws.Range("H8").Value = 25.623216
$ws.Range("I8").Value = 0.2580047039168495
$ws.Range("J8").Value = 0.2580047039168495
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.198887402592
$ws.Range("R8").Value = 1.789986623328
$ws.Range("S8").Value = 0.002402207099563116
$ws.Range("T8").Value = 0.002402207099563116
$ws.Range("G9").Value = 8.541072
$ws.Range("H9").Value = 25.623216
$ws.Range("I9").Value = 0.2580047039168495
$ws.Range("J9").Value = 0.2580047039168495
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 1.135273596192
$ws.Range("R9").Value = 10.217462365728
$ws.Range("S9").Value = 0.01371209165174481
$ws.Range("T9").Value = 0.01371209165174481
$ws.Range("G10").Value = 8.541072
$ws.Range("H10").Value = 25.623216
$ws.Range("I10").Value = 0.2580047039168495
$ws.Range("J10").Value = 0.2580047039168495
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 20.026980356544
$ws.Range("R10").Value = 180.242823208896
$ws.Range("S10").Value = 0.2418904051655416
$ws.Range("T10").Value = 0.2418904051655416
